# Updates cryptocurrency price/volume data in the "cryptos" worksheet
# to reflect the latest GitHub Actions scrape (commit: "Updated cryptos
# list on Sat Apr 22 15:28:35 UTC 2023 with GitHub Actions").
#
# Column D (Price) values are formatted as plain text that often look
# like numbers (e.g. "27.469.49", "0.06574", "47.50"); NumberFormat is
# forced to Text ("@") before assignment so Excel does not silently
# coerce them into numeric values and strip significant trailing
# zeros / thousands separators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.469.49'
$ws.Range("E2").Value = '  -2.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.864.91'
$ws.Range("E3").Value = '  -2.51%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.53'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4713'
$ws.Range("E7").Value = '  +1.54%  '
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.50'
$ws.Range("E9").Value = '  -10.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08024'
$ws.Range("E10").Value = '  -4.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.020'
$ws.Range("E11").Value = '  -2.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.64'
$ws.Range("E12").Value = '  -1.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.890.58'
$ws.Range("E13").Value = '  -1.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.957'
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.194'
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.75'
$ws.Range("E17").Value = '  -3.04%  '
$ws.Range("E18").Value = '  -2.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06574'
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.29'
$ws.Range("E20").Value = '  -3.26%  '
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.514'
$ws.Range("E22").Value = '  -3.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.483.86'
$ws.Range("E23").Value = '  -2.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.98'
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.301'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.078.01'
$ws.Range("E26").Value = '  -2.93%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.99'
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.26'
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.082'
$ws.Range("E29").Value = '  -2.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.547'
$ws.Range("E30").Value = '  -3.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.36'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("E32").Value = '  -1.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9568'
$ws.Range("E33").Value = '  -1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.468'
$ws.Range("E34").Value = '  +1.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.595'
$ws.Range("E35").Value = '  -1.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.297'
$ws.Range("E36").Value = '  -4.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06100'
$ws.Range("E37").Value = '  -0.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02248'
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.220'
$ws.Range("E39").Value = '  -4.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.099'
$ws.Range("E40").Value = '  -8.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1897'
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.33'
$ws.Range("E44").Value = '  -5.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.266'
$ws.Range("E45").Value = '  -2.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5674'
$ws.Range("E46").Value = '  -3.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.21'
$ws.Range("E47").Value = '  -4.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.416'
$ws.Range("E48").Value = '  -0.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.945'
$ws.Range("E49").Value = '  -3.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06780'
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("E51").Value = '  -1.21%  '